$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1151.1333
$ws.Range("J17").Value = 1151.1333
$ws.Range("L17").Value = 3453.3999
$ws.Range("N17").Value = -3789.3999
$ws.Range("H19").Value = 1062.375
$ws.Range("I19").Value = 1025
$ws.Range("J19").Value = 1174.5
$ws.Range("K19").Value = 1025
$ws.Range("L19").Value = 1174.5
$ws.Range("M19").Value = -850
$ws.Range("N19").Value = -1524.5
$ws.Range("H51").Value = 7999.4
$ws.Range("I51").Value = 7999.4
$ws.Range("K51").Value = 7999.4
$ws.Range("M51").Value = -7515.4
$ws.Range("H62").Value = 2955.25
$ws.Range("I62").Value = 2943.6667
$ws.Range("J62").Value = 2990
$ws.Range("K62").Value = 2943.6667
$ws.Range("L62").Value = 2990
$ws.Range("M62").Value = -2319.6667
$ws.Range("N62").Value = -4238
$ws.Range("H64").Value = 5000
$ws.Range("I64").Value = 5000
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 5000
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -4752
$ws.Range("N64").Value = -5496
$ws.Range("H65").Value = 2955.25
$ws.Range("I65").Value = 2943.6667
$ws.Range("J65").Value = 2990
$ws.Range("K65").Value = 14718.3335
$ws.Range("L65").Value = 14950
$ws.Range("M65").Value = -11598.3335
$ws.Range("N65").Value = -21190
$ws.Range("H67").Value = 5000
$ws.Range("I67").Value = 5000
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 5000
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -4142
$ws.Range("N67").Value = -6716
$ws.Range("H74").Value = 102781.6
$ws.Range("I74").Value = 3477
$ws.Range("J74").Value = 168984.67
$ws.Range("K74").Value = 3477
$ws.Range("L74").Value = 168984.67
$ws.Range("M74").Value = -2541
$ws.Range("N74").Value = -170856.67
$ws.Range("H77").Value = 102781.6
$ws.Range("I77").Value = 3477
$ws.Range("J77").Value = 168984.67
$ws.Range("K77").Value = 17385
$ws.Range("L77").Value = 844923.3500000001
$ws.Range("M77").Value = -12705
$ws.Range("N77").Value = -854283.3500000001
$ws.Range("H132").Value = 1726.7778
$ws.Range("I132").Value = 1518.2609
$ws.Range("J132").Value = 2925.75
$ws.Range("K132").Value = 4554.7827
$ws.Range("L132").Value = 8777.25
$ws.Range("M132").Value = -2024.7827
$ws.Range("N132").Value = -13837.25
$ws.Range("H137").Value = 3367
$ws.Range("J137").Value = 3284.7856
$ws.Range("L137").Value = 9854.356800000001
$ws.Range("N137").Value = -14954.3568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 500
$ws.Range("I5").Value = 500
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 500
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -388
$ws.Range("N5").Value = ""
$ws.Range("H74").Value = 1149.6666
$ws.Range("I74").Value = 1224.5
$ws.Range("K74").Value = 1224.5
$ws.Range("M74").Value = -350.5
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").Value = ""
$ws.Range("H77").Value = 1149.6666
$ws.Range("I77").Value = 1224.5
$ws.Range("K77").Value = 6122.5
$ws.Range("M77").Value = -1754.5
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").Value = ""
$ws.Range("H139").Value = 29000
$ws.Range("J139").Value = 29000
$ws.Range("L139").Value = 29000
$ws.Range("N139").Value = -39280

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 500
$ws.Range("I4").Value = 500
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 500
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -385
$ws.Range("N4").Value = ""
$ws.Range("H134").Value = 3820.6843
$ws.Range("I134").Value = 3310.7778
$ws.Range("K134").Value = 9932.3334
$ws.Range("M134").Value = -7397.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 42.444443
$ws.Range("J7").Value = 75
$ws.Range("L7").Value = 75
$ws.Range("N7").Value = -301
$ws.Range("H22").Value = 639.2
$ws.Range("I22").Value = 460.2
$ws.Range("J22").Value = 818.2
$ws.Range("K22").Value = 460.2
$ws.Range("L22").Value = 818.2
$ws.Range("M22").Value = -110.2
$ws.Range("N22").Value = -1518.2
$ws.Range("H32").Value = 2372.5
$ws.Range("I32").Value = 1666.3334
$ws.Range("K32").Value = 1666.3334
$ws.Range("M32").Value = -1350.3334
$ws.Range("H33").Value = 575
$ws.Range("I33").Value = 575
$ws.Range("K33").Value = 575
$ws.Range("M33").Value = -196
$ws.Range("H42").Value = 4000
$ws.Range("I42").Value = 3000
$ws.Range("J42").Value = 5000
$ws.Range("K42").Value = 3000
$ws.Range("L42").Value = 5000
$ws.Range("M42").Value = -2407
$ws.Range("N42").Value = -6186
$ws.Range("H44").Value = 5971
$ws.Range("J44").Value = 5971
$ws.Range("L44").Value = 5971
$ws.Range("N44").Value = -6855
$ws.Range("H45").Value = 43437
$ws.Range("J45").Value = 43437
$ws.Range("L45").Value = 43437
$ws.Range("N45").Value = -44623
$ws.Range("H54").Value = 8300
$ws.Range("I54").Value = 8300
$ws.Range("K54").Value = 8300
$ws.Range("M54").Value = -7642
$ws.Range("H55").Value = 4100
$ws.Range("I55").Value = 4100
$ws.Range("K55").Value = 4100
$ws.Range("M55").Value = -3785
$ws.Range("H57").Value = 25000
$ws.Range("J57").Value = 25000
$ws.Range("L57").Value = 25000
$ws.Range("N57").Value = -26120
$ws.Range("H63").Value = 40000
$ws.Range("J63").Value = 40000
$ws.Range("L63").Value = 40000
$ws.Range("N63").Value = -41372
$ws.Range("H66").Value = 40000
$ws.Range("J66").Value = 40000
$ws.Range("L66").Value = 120000
$ws.Range("N66").Value = -126864
$ws.Range("H81").Value = 100777
$ws.Range("J81").Value = 100777
$ws.Range("L81").Value = 100777
$ws.Range("N81").Value = -102773
$ws.Range("H84").Value = 100777
$ws.Range("J84").Value = 100777
$ws.Range("L84").Value = 302331
$ws.Range("N84").Value = -312315
$ws.Range("H88").Value = 46856
$ws.Range("J88").Value = 46856
$ws.Range("L88").Value = 46856
$ws.Range("N88").Value = -47668
$ws.Range("H91").Value = 46856
$ws.Range("J91").Value = 46856
$ws.Range("L91").Value = 46856
$ws.Range("N91").Value = -49664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41.363636
$ws.Range("I12").Value = 18.833334
$ws.Range("J12").Value = 68.40000000000001
$ws.Range("K12").Value = 56.500002
$ws.Range("L12").Value = 205.2
$ws.Range("M12").Value = 116.499998
$ws.Range("N12").Value = -551.2
$ws.Range("H98").Value = 274.5
$ws.Range("J98").Value = 274.5
$ws.Range("L98").Value = 823.5
$ws.Range("N98").Value = -3819.5
$ws.Range("H120").Value = 999
$ws.Range("I120").Value = 999
$ws.Range("K120").Value = 2997
$ws.Range("M120").Value = 1841
$ws.Range("H136").Value = 3525
$ws.Range("I136").Value = 3050
$ws.Range("J136").Value = 4950
$ws.Range("K136").Value = 9150
$ws.Range("L136").Value = 14850
$ws.Range("M136").Value = -4050
$ws.Range("N136").Value = -25050
$ws.Range("H140").Value = 2045.7142
$ws.Range("J140").Value = 2950
$ws.Range("L140").Value = 8850
$ws.Range("N140").Value = -19210

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 504
$ws.Range("I5").Value = 504
$ws.Range("K5").Value = 504
$ws.Range("M5").Value = -392

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1049.8
$ws.Range("I22").Value = 855.75
$ws.Range("K22").Value = 855.75
$ws.Range("M22").Value = -560.75
$ws.Range("H27").Value = 1049.8
$ws.Range("I27").Value = 855.75
$ws.Range("K27").Value = 855.75
$ws.Range("M27").Value = -748.75
$ws.Range("H46").Value = 4501
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 4501
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 4501
$ws.Range("M46").Value = ""
$ws.Range("N46").Value = -4877
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H132").Value = 7694
$ws.Range("I132").Value = 4999.6665
$ws.Range("K132").Value = 14998.9995
$ws.Range("M132").Value = -12468.9995
